$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend row 3 with new sentiment "up/down" info columns ---
$ws.Range("X3").Value = -0.93999500000001035
$ws.Range("Y3").Value = "Down"

# --- Add new row 4: latest trading day data (Long/Hold trade) ---

# A4 needs the same date style as A2/A3 - copy formatting from A3 first
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 42633.890567129631

$ws.Range("B4").Value = 19
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Random"
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.86

# S4 needs the same percentage style as S2/S3 - copy formatting from S3 first
$ws.Range("S3").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 0.0262

$ws.Range("T4").Value = -2.66
$ws.Range("U4").Value = 15.05
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0
